# Data cleaning: rename the worksheet tab from "query" to "Sheet1" and
# leave the selection parked on D30 (matches the saved cursor position
# after the cleanup pass).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Renaming the sheet also updates the hidden "query" defined name's
# reference (query!$A$1:$Y$360 -> Sheet1!$A$1:$Y$360) automatically.
$ws.Name = "Sheet1"

# Move/record the active selection at D30.
$ws.Range("D30").Select()
